$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: find the first occurrence of literal $text in the whole document
# and return the Range for it (or $null if not found).
# ---------------------------------------------------------------------------
function Find-First($text) {
    $r = $d.Content
    $ok = $r.Find.Execute($text, $false, $false, $false, $false, $false, $true, 1, $false, $null, 0)
    if (-not $ok) { return $null }
    return $r
}

# ---------------------------------------------------------------------------
# Helper: insert $newText at the (collapsed) insertion point $ip, giving the
# inserted run the exact same run-formatting (rPr) as $donor (another Range,
# typically a single already-existing run elsewhere in the document whose
# formatting we want to clone exactly, incl. rFonts eastAsia/cs that plain
# Font.* property writes do not reproduce faithfully).
#
# Returns the Range that now covers the freshly-inserted text.
# ---------------------------------------------------------------------------
function Insert-WithFormattingOf($ip, $donor, $newText) {
    $start = $ip.Start
    # Clone formatting+text of the donor run at the insertion point; this
    # inserts (does not overwrite neighbouring text) because $ip is a
    # zero-length range.
    $ip.FormattedText = $donor.FormattedText
    $donorLen = $donor.End - $donor.Start
    $placed = $d.Range($start, $start + $donorLen)
    # Now stamp the text we actually want onto the freshly inserted run
    # while preserving the run's formatting.
    $placed.Text = $newText
    return $d.Range($start, $start + $newText.Length)
}

# Donor ranges used purely as formatting templates (their own text is
# irrelevant - Insert-WithFormattingOf overwrites it immediately after).
$donorBlue = Find-First("<al>")
$donorGray = Find-First("amp;")

# ===========================================================================
# Change 1: "...Jay veu que celuy qui maprist a les mouler" ->
#           "...Jay veu que <pn>celuy qui maprist a les mouler</pn>"
# ===========================================================================
$r1 = Find-First("en forme de S Jay veu que celuy qui maprist a les mouler")
$splitAt = $r1.Start + ("en forme de S Jay veu que ").Length
$ip1 = $d.Range($splitAt, $splitAt)
$tag1 = Insert-WithFormattingOf $ip1 $donorBlue "<pn>"

$endAt = $tag1.End + ("celuy qui maprist a les mouler").Length
$ip2 = $d.Range($endAt, $endAt)
Insert-WithFormattingOf $ip2 $donorBlue "</pn>" | Out-Null

# ===========================================================================
# Change 2: "<m>terre</m> fresche &amp; humide</env>" ->
#           "<m>terre fresche &amp; humide</m></env>"
# (delete the "</m>" run - its neighbours merge - then rebuild it in front
#  of "</env>")
# ===========================================================================
$r2 = Find-First("</m> fresche")
$delRange2 = $d.Range($r2.Start, $r2.Start + 4)   # "</m>"
$delRange2.Text = ""

$envSearchStart = $delRange2.End
$r2env = $d.Range($envSearchStart, $d.Content.End)
$r2env.Find.Execute("</env>", $false, $false, $false, $false, $false, $true, 1, $false, $null, 0) | Out-Null
$ip3 = $d.Range($r2env.Start, $r2env.Start)
Insert-WithFormattingOf $ip3 $donorBlue "</m>" | Out-Null

# ===========================================================================
# Change 3: "quun quart" -> "quun <ms><tmp>quart"
# ===========================================================================
$r3 = Find-First("quun quart")
$splitAt3 = $r3.Start + ("quun ").Length
$ip4 = $d.Range($splitAt3, $splitAt3)
Insert-WithFormattingOf $ip4 $donorBlue "<ms><tmp>" | Out-Null

# ===========================================================================
# Change 4: "d<ms>heure</ms>" -> "dheure</tmp></ms>"
# (delete the "<ms>" run - its neighbours merge - then rebuild the closing
#  tags in front of "</ms>")
# ===========================================================================
$r4 = Find-First("d<ms>heure</ms>")
$msStart = $r4.Start + ("d").Length
$delRange4 = $d.Range($msStart, $msStart + ("<ms>").Length)
$delRange4.Text = ""

$msEnvSearchStart = $delRange4.End
$r4ms = $d.Range($msEnvSearchStart, $d.Content.End)
$r4ms.Find.Execute("</ms>", $false, $false, $false, $false, $false, $true, 1, $false, $null, 0) | Out-Null
$ip5 = $d.Range($r4ms.Start, $r4ms.Start)
Insert-WithFormattingOf $ip5 $donorBlue "</tmp>" | Out-Null

Write-Output "Done"
